# Actualización automática del mapa (2025-08-25 13:59:21)
# Adds a new record row (row 89) to the single "AYKO" worksheet,
# mirroring the structure of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

# Columns A-H, J-L, O-P hold text values in this sheet (even when they
# look numeric, e.g. case numbers, dates-as-text, phone/OT numbers).
# Force a text number format before assigning so Excel doesn't
# auto-convert these into numbers or dates.
$ws.Range("A$row`:H$row").NumberFormat = "@"
$ws.Range("J$row`:L$row").NumberFormat = "@"
$ws.Range("O$row`:P$row").NumberFormat = "@"

$ws.Range("A$row").Value = "-563"
$ws.Range("B$row").Value = "8/25/2025"
$ws.Range("C$row").Value = "Av Castañares 4520"
$ws.Range("D$row").Value = "8"
$ws.Range("E$row").Value = "809157022"
$ws.Range("F$row").Value = "AYKO"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Cambiar pasante"
$ws.Range("I$row").Value = 1
$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Pasante"
$ws.Range("M$row").Value = -58.470249
$ws.Range("N$row").Value = -34.664835
$ws.Range("O$row").Value = "Boedo"
$ws.Range("P$row").Value = "Capital Sur"

# Reset the style of the new row back to the default (Normal) so no
# extra/residual number-formatting is left behind, matching the look
# of the other data rows (which carry no explicit style index).
$ws.Range("A$row`:P$row").Style = "Normal"
